# Update countries & provincias Spain
# Refresh the COVID country table (rows keep their row number / rank in
# column A is untouched since that's just an ordinal; the country name in
# column A and the stats in B:H are refreshed to the new snapshot). Where a
# country's total cases overtook its neighbour, the two countries swap rows
# so the sheet stays sorted by "Casos totales" (column B) descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow {
    param(
        [int]$Row,
        [string]$Country,
        [double]$CasosTotales,
        [double]$NuevosCasos,
        [double]$CasosActivos,
        [double]$Recuperados,
        [double]$MuertesHoy,
        [double]$Muertes
    )

    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = 0
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# row, country, Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes
Set-CountryRow 4   "Estados Unidos"         4635226 241  2285613 2194307 21  155306
Set-CountryRow 13  "Iran"                    304204 2674  263519   23919 197  16766
Set-CountryRow 14  "Reino Unido"             302301    0       0       0   0  45999
Set-CountryRow 43  "Emiratos Arabes Unidos"   60506  283   53909    6246   2    351
Set-CountryRow 46  "Rumania"                  50886 1295   27007   21536  39   2343
Set-CountryRow 47  "Portugal"                 50868    0   36140   13001   0   1727
Set-CountryRow 54  "Afganistan"               36675  133   25509    9894   1   1272
Set-CountryRow 56  "Suiza"                    35232  210   31100    2152   0   1980
Set-CountryRow 57  "Ghana"                    35142    0   31286    3681   0    175
Set-CountryRow 73  "El Salvador"              16632  402    8362    7822   9    448
Set-CountryRow 74  "Etiopia"                  16615    0    6763    9589   0    263
Set-CountryRow 75  "Chequia"                  16371    0   11482    4510   0    379
Set-CountryRow 125 "Eslovenia"                 2156   17    1797     241   1    118

# Row 162 (Vietnam): only Recuperados / Muertes hoy / Muertes changed.
$ws.Cells.Item(162, 5).Value = 135
$ws.Cells.Item(162, 7).Value = 1
$ws.Cells.Item(162, 8).Value = 1

# Footer timestamp string (cell A1) advances from 11:18 to 12:35.
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 12:35"
